# CASEFLOW-4586 Update test spreadsheets with White River
#
# 1) "RO Non-Availability Dates" sheet: insert a new column F for
#    RO05 / White River Junction, VT (shifts existing columns F.. right).
# 2) "RO Allocations" sheet: insert a new row 7 for
#    RO05 / White River Junction, VT (shifts existing rows 7.. down).
# 3) Re-point the active sheet/selection to match the authored session.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("RO Non-Availability Dates")
$ws3 = $wb.Worksheets.Item("RO Allocations")

# --- Sheet 1: insert column F (RO05) ----------------------------------
$ws1.Columns.Item(6).Insert()
$ws1.Columns.Item(6).ColumnWidth = 14.67
$ws1.Columns.Item(7).ColumnWidth = 14.67

# New column F inherits column G's (old F's) formatting in row 3 so the
# header cell styling lines up with its neighbour.
$ws1.Cells.Item(3, 7).Copy()
$ws1.Cells.Item(3, 6).PasteSpecial(-4122)

# Values - set row 3 (city) before row 2 (RO code) so the shared-string
# table gets "White River Junction, VT" before "RO05".
$ws1.Cells.Item(3, 6).Value = "White River Junction, VT"
$ws1.Cells.Item(2, 6).Value = "RO05"

# --- Sheet 3: insert row 7 (RO05) --------------------------------------
$ws3.Rows.Item(7).Insert()

# New row 7 inherits row 8's (old row 7's) formatting across the used
# columns (A:N) only, so we don't drag formatting across the whole row.
$ws3.Range("A8:N8").Copy()
$ws3.Range("A7:N7").PasteSpecial(-4122)

$ws3.Range("B7").Value = "White River Junction, VT"
$ws3.Range("C7").Value = "RO05"
$ws3.Range("D7").Value = 0
$ws3.Range("E7").Value = 0
$ws3.Range("F7").Value = 0
$ws3.Range("G7").Value = 0
$ws3.Range("H7").Value = 0

# --- Active sheet / selection ------------------------------------------
$ws3.Activate()
$ws3.Range("E12").Select()

$ws1.Activate()
$ws1.Range("G12").Select()
